# feat: add missing strings
#
# Inserts 2 new i18n rows (appMenu.frame, appMenu.feedback) right after the
# existing "appMenu.bookmark" row and before "appMenu.bookmark.undo", and
# appends 3 new i18n rows (app.pornhub.fontSize, app.pornhub.vertical,
# app.pornhub.colorRevert) at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make room for two new rows right after row 25 (appMenu.bookmark) ---
# Shift the existing rows 26:57 down to 28:59, carrying their formatting
# (style + row height) with them, then overwrite 26:27 with the new data.
$ws.Range("A26:E57").Copy($ws.Range("A28:E59"))

$ws.Range("A26").Value = "appMenu.frame"
$ws.Range("B26").Value = "框架引用"
$ws.Range("C26").Value = "Use in your website"

$ws.Range("A27").Value = "appMenu.feedback"
$ws.Range("B27").Value = "反馈"
$ws.Range("C27").Value = "Feedback"

# --- 2) Append three new rows at the bottom of the table (now rows 60-62) ---
# Row 59 ("app.urlcleaner.ruleTitle" / Rules) is now the last data row;
# clone its formatting down for the three new rows.
$ws.Range("A59:E59").Copy($ws.Range("A60:E62"))

$ws.Range("A60").Value = "app.pornhub.fontSize"
$ws.Range("B60").Value = "字体大小"
$ws.Range("C60").Value = "Font Size"

$ws.Range("A61").Value = "app.pornhub.vertical"
$ws.Range("B61").Value = "竖直排列"
$ws.Range("C61").Value = "Vertical Layout"

$ws.Range("A62").Value = "app.pornhub.colorRevert"
$ws.Range("B62").Value = "颜色反转"
$ws.Range("C62").Value = "Color Revert"

# --- 3) Restore the uniform row height on every row touched by the block
# copies above (a multi-row Range.Copy doesn't always carry the source
# row's `ht` over to new destination rows), matching the rest of the table.
$ws.Range("A26:E62").RowHeight = 20.1
